$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2476489028213166
$ws.Range("C2").Value = 0.4576802507836991
$ws.Range("J2").Value = 0.02821316614420063
$ws.Range("O2").Value = 0.003134796238244514
$ws.Range("P2").Value = 0.1755485893416928
$ws.Range("S2").Value = 0.0877742946708464

# Row 3
$ws.Range("B3").Value = 0.02580645161290323
$ws.Range("C3").Value = 0.03870967741935484
$ws.Range("J3").Value = 0.03870967741935484
$ws.Range("P3").Value = 0.7612903225806451
$ws.Range("S3").Value = 0.1354838709677419

# Row 4
$ws.Range("J4").Value = 0.02173913043478261
$ws.Range("P4").Value = 0.7826086956521739
$ws.Range("S4").Value = 0.1956521739130435

# Row 6
$ws.Range("B6").Value = 0.05882352941176471
$ws.Range("D6").Value = 0.009049773755656109
$ws.Range("F6").Value = 0.08597285067873303
$ws.Range("J6").Value = 0.2443438914027149
$ws.Range("O6").Value = 0.04072398190045249
$ws.Range("Q6").Value = 0.1176470588235294
$ws.Range("R6").Value = 0.09954751131221719
$ws.Range("S6").Value = 0.3438914027149321

# Row 7
$ws.Range("B7").Value = 0.07017543859649122
$ws.Range("D7").Value = 0.02339181286549707
$ws.Range("E7").Value = 0.005847953216374269
$ws.Range("F7").Value = 0.07017543859649122
$ws.Range("J7").Value = 0.09941520467836257
$ws.Range("O7").Value = 0.02339181286549707
$ws.Range("Q7").Value = 0.1754385964912281
$ws.Range("R7").Value = 0.08187134502923976
$ws.Range("S7").Value = 0.4502923976608187

# Row 8
$ws.Range("B8").Value = 0.09266409266409266
$ws.Range("D8").Value = 0.01544401544401544
$ws.Range("E8").Value = 0.001930501930501931
$ws.Range("F8").Value = 0.05598455598455598
$ws.Range("J8").Value = 0.1042471042471042
$ws.Range("O8").Value = 0.02316602316602316
$ws.Range("Q8").Value = 0.1428571428571428
$ws.Range("R8").Value = 0.1138996138996139
$ws.Range("S8").Value = 0.4498069498069498

# Row 9
$ws.Range("B9").Value = 0.155688622754491
$ws.Range("D9").Value = 0.01796407185628742
$ws.Range("F9").Value = 0.08383233532934131
$ws.Range("J9").Value = 0.0718562874251497
$ws.Range("O9").Value = 0.005988023952095809
$ws.Range("Q9").Value = 0.1377245508982036
$ws.Range("R9").Value = 0.1017964071856287
$ws.Range("S9").Value = 0.4251497005988024

# Row 10
$ws.Range("B10").Value = 0.109715242881072
$ws.Range("D10").Value = 0.02596314907872697
$ws.Range("E10").Value = 0.001675041876046901
$ws.Range("F10").Value = 0.05946398659966499
$ws.Range("J10").Value = 0.1080402010050251
$ws.Range("O10").Value = 0.01675041876046901
$ws.Range("Q10").Value = 0.2093802345058627
$ws.Range("R10").Value = 0.09631490787269682
$ws.Range("S10").Value = 0.3726968174204355

# Row 11
$ws.Range("G11").Value = 0.1346938775510204
$ws.Range("J11").Value = 0.08571428571428572
$ws.Range("K11").Value = 0.1795918367346939
$ws.Range("L11").Value = 0.5673469387755102
$ws.Range("S11").Value = 0.0326530612244898

# Row 12
$ws.Range("G12").Value = 0.7534246575342466
$ws.Range("J12").Value = 0.1780821917808219
$ws.Range("L12").Value = 0.04794520547945205
$ws.Range("S12").Value = 0.02054794520547945

# Row 13
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2156862745098039
$ws.Range("S13").Value = 0.1176470588235294

# Row 15
$ws.Range("F15").Value = 0.02666666666666667
$ws.Range("H15").Value = 0.1377777777777778
$ws.Range("I15").Value = 0.04444444444444445
$ws.Range("J15").Value = 0.3422222222222222
$ws.Range("K15").Value = 0.06222222222222222
$ws.Range("M15").Value = 0.01333333333333333
$ws.Range("N15").Value = 0.004444444444444444
$ws.Range("O15").Value = 0.05777777777777778
$ws.Range("S15").Value = 0.3111111111111111

# Row 16
$ws.Range("F16").Value = 0.02403846153846154
$ws.Range("H16").Value = 0.25
$ws.Range("I16").Value = 0.08173076923076923
$ws.Range("J16").Value = 0.2884615384615384
$ws.Range("K16").Value = 0.09134615384615384
$ws.Range("M16").Value = 0.01923076923076923
$ws.Range("N16").Value = 0.004807692307692308
$ws.Range("O16").Value = 0.0673076923076923
$ws.Range("S16").Value = 0.1730769230769231

# Row 17
$ws.Range("F17").Value = 0.01951219512195122
$ws.Range("H17").Value = 0.2097560975609756
$ws.Range("I17").Value = 0.06585365853658537
$ws.Range("J17").Value = 0.4121951219512195
$ws.Range("K17").Value = 0.08536585365853659
$ws.Range("M17").Value = 0.02439024390243903
$ws.Range("N17").Value = 0.004878048780487805
$ws.Range("O17").Value = 0.08048780487804878
$ws.Range("S17").Value = 0.0975609756097561

# Row 18
$ws.Range("F18").Value = 0.008771929824561403
$ws.Range("H18").Value = 0.2368421052631579
$ws.Range("I18").Value = 0.09649122807017543
$ws.Range("J18").Value = 0.4078947368421053
$ws.Range("K18").Value = 0.06140350877192982
$ws.Range("M18").Value = 0.01754385964912281
$ws.Range("O18").Value = 0.05263157894736842
$ws.Range("S18").Value = 0.1184210526315789

# Row 19
$ws.Range("F19").Value = 0.021671826625387
$ws.Range("H19").Value = 0.2306501547987616
$ws.Range("I19").Value = 0.07043343653250773
$ws.Range("J19").Value = 0.3653250773993808
$ws.Range("K19").Value = 0.08823529411764706
$ws.Range("M19").Value = 0.02321981424148607
$ws.Range("O19").Value = 0.06578947368421052
$ws.Range("S19").Value = 0.1346749226006192
